# Refresh the crypto Price (D) / Volume(1h) (E) columns with the latest scrape.
# Commit: "Updated cryptos list on Sun Dec 17 06:46:51 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $range = $ws.Range($Cell)
    if ($Text -match "^[+-]?[0-9]*\.?[0-9]+$") {
        # Looks like a plain number (e.g. "240.36") - without help Excel would
        # silently coerce the assignment to a numeric value (and drop formatting
        # such as trailing zeros, e.g. "6.10" -> 6.1). Force text storage, then
        # restore the default "Normal" style so no formatting change is introduced.
        $range.NumberFormat = "@"
        $range.Value = $Text
        $range.Style = "Normal"
    } else {
        # Already safe to assign as-is (contains extra separators, a "%" sign, etc.)
        $range.Value = $Text
    }
}

Set-TextValue "D2" "41.942.47"
Set-TextValue "E2" "  -0.72%  "
Set-TextValue "D3" "2.205.79"
Set-TextValue "E3" "  -1.63%  "
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "240.36"
Set-TextValue "E5" "  -2.43%  "
Set-TextValue "E6" "  -0.80%  "
Set-TextValue "D7" "72.87"
Set-TextValue "E7" "  -1.82%  "
Set-TextValue "E8" "  +0.18%  "
Set-TextValue "E9" "  -1.81%  "
Set-TextValue "D10" "42.84"
Set-TextValue "E10" "  +2.23%  "
Set-TextValue "D11" "0.0947"
Set-TextValue "E11" "  +0.81%  "
Set-TextValue "D12" "7.09"
Set-TextValue "E12" "  -0.79%  "
Set-TextValue "E13" "  -0.47%  "
Set-TextValue "D14" "2.537.79"
Set-TextValue "E14" "  -1.51%  "
Set-TextValue "D15" "14.17"
Set-TextValue "E15" "  -1.95%  "
Set-TextValue "D16" "0.837"
Set-TextValue "E16" "  -1.52%  "
Set-TextValue "D17" "2.211.61"
Set-TextValue "E17" "  -0.91%  "
Set-TextValue "D18" "41.813.23"
Set-TextValue "E18" "  -0.53%  "
Set-TextValue "E19" "  +9.27%  "
Set-TextValue "D20" "72.52"
Set-TextValue "E20" "  +0.82%  "
Set-TextValue "D21" "6.10"
Set-TextValue "E21" "  -0.66%  "
Set-TextValue "D22" "10.24"
Set-TextValue "E22" "  +17.55%  "
Set-TextValue "D23" "228.55"
Set-TextValue "E23" "  -1.43%  "
Set-TextValue "D24" "2.06"
Set-TextValue "E24" "  -7.62%  "
Set-TextValue "E25" "  +0.16%  "
Set-TextValue "D26" "11.45"
Set-TextValue "E26" "  +1.18%  "
Set-TextValue "E27" "  +0.09%  "
Set-TextValue "D28" "2.25"
Set-TextValue "E28" "  -2.58%  "
Set-TextValue "E29" "  +0.61%  "
Set-TextValue "D30" "167.05"
Set-TextValue "E30" "  -1.25%  "
Set-TextValue "D31" "20.49"
Set-TextValue "E31" "  -0.61%  "
Set-TextValue "D32" "5.54"
Set-TextValue "D33" "0.0788"
Set-TextValue "E33" "  -3.92%  "
Set-TextValue "E34" "  -0.39%  "
Set-TextValue "D35" "28.91"
Set-TextValue "E35" "  -5.09%  "
Set-TextValue "E36" "  -7.32%  "
Set-TextValue "E37" "  -5.70%  "
Set-TextValue "D38" "0.0299"
Set-TextValue "E38" "  -1.79%  "
Set-TextValue "D39" "12.85"
Set-TextValue "E39" "  -7.02%  "
Set-TextValue "D40" "65.16"
Set-TextValue "E40" "  +5.09%  "
Set-TextValue "E41" "  -3.40%  "
Set-TextValue "D42" "5.60"
Set-TextValue "E42" "  -3.19%  "
Set-TextValue "E43" "  -1.88%  "
Set-TextValue "D44" "8.67"
Set-TextValue "E44" "  +0.53%  "
Set-TextValue "D45" "103.65"
Set-TextValue "E45" "  -3.33%  "
Set-TextValue "E46" "  -2.00%  "
Set-TextValue "D47" "2.38"
Set-TextValue "E47" "  +4.61%  "
Set-TextValue "E48" "  -1.25%  "
Set-TextValue "E49" "  -1.02%  "
Set-TextValue "E50" "  +0.16%  "
Set-TextValue "D51" "2.413.35"
Set-TextValue "E51" "  -1.53%  "
